# Auto-generated Excel COM-interop script applying the Ravana_Profits edits
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1855.3334
$ws.Range("I28").Value = 1796.5
$ws.Range("K28").Value = 1796.5
$ws.Range("M28").Value = -1311.5
$ws.Range("H43").Value = 6661
$ws.Range("I43").Value = 6661
$ws.Range("K43").Value = 6661
$ws.Range("M43").Value = -6592
$ws.Range("H58").Value = 7743.2
$ws.Range("J58").Value = 11972
$ws.Range("L58").Value = 35916
$ws.Range("N58").Value = -36216
$ws.Range("H98").Value = 1199.6
$ws.Range("I98").Value = 1332.6666
$ws.Range("K98").Value = 1332.6666
$ws.Range("M98").Value = 165.3334
$ws.Range("H112").Value = 1646.2273
$ws.Range("J112").Value = 1785.1052
$ws.Range("L112").Value = 5355.3156
$ws.Range("N112").Value = -7571.3156
$ws.Range("H122").Value = 1199.6
$ws.Range("I122").Value = 1332.6666
$ws.Range("K122").Value = 3997.9998
$ws.Range("M122").Value = -1547.9998
$ws.Range("H132").Value = 1188.7222
$ws.Range("I132").Value = 1188.7222
$ws.Range("K132").Value = 3566.1666
$ws.Range("M132").Value = -1036.1666
$ws.Range("H138").Value = 2841.3281
$ws.Range("I138").Value = 2064.5833
$ws.Range("J138").Value = 3020.577
$ws.Range("K138").Value = 6193.749899999999
$ws.Range("L138").Value = 9061.731
$ws.Range("M138").Value = -1053.749899999999
$ws.Range("N138").Value = -19341.731

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1615
$ws.Range("I2").Value = 1810.091
$ws.Range("K2").Value = 1810.091
$ws.Range("M2").Value = -1697.091
$ws.Range("H32").Value = 4804.5
$ws.Range("I32").Value = 3416.8684
$ws.Range("K32").Value = 3416.8684
$ws.Range("M32").Value = -3129.8684
$ws.Range("H74").Value = 4342.3335
$ws.Range("I74").Value = 4013
$ws.Range("K74").Value = 4013
$ws.Range("M74").Value = -3139
$ws.Range("H77").Value = 4342.3335
$ws.Range("I77").Value = 4013
$ws.Range("K77").Value = 20065
$ws.Range("M77").Value = -15697
$ws.Range("H102").Value = 1418.4445
$ws.Range("I102").Value = 1418.4445
$ws.Range("K102").Value = 1418.4445
$ws.Range("M102").Value = 203.5554999999999
$ws.Range("H110").Value = 1199.5
$ws.Range("I110").Value = 1199.5
$ws.Range("K110").Value = 1199.5
$ws.Range("M110").Value = 845.5
$ws.Range("H116").Value = 1615
$ws.Range("I116").Value = 1810.091
$ws.Range("K116").Value = 1810.091
$ws.Range("M116").Value = 483.9090000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1615
$ws.Range("I3").Value = 1810.091
$ws.Range("K3").Value = 1810.091
$ws.Range("M3").Value = -1696.091
$ws.Range("H134").Value = 2258.3462
$ws.Range("I134").Value = 2066.15
$ws.Range("K134").Value = 6198.450000000001
$ws.Range("M134").Value = -3663.450000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1613.7142
$ws.Range("J31").Value = 1699.25
$ws.Range("L31").Value = 1699.25
$ws.Range("N31").Value = -2289.25
$ws.Range("H34").Value = 1613.7142
$ws.Range("J34").Value = 1699.25
$ws.Range("L34").Value = 1699.25
$ws.Range("N34").Value = -2103.25
$ws.Range("H86").Value = 9641.308000000001
$ws.Range("I86").Value = 11999.571
$ws.Range("J86").Value = 6890
$ws.Range("K86").Value = 11999.571
$ws.Range("L86").Value = 6890
$ws.Range("M86").Value = -10876.571
$ws.Range("N86").Value = -9136
$ws.Range("H89").Value = 9641.308000000001
$ws.Range("I89").Value = 11999.571
$ws.Range("J89").Value = 6890
$ws.Range("K89").Value = 59997.855
$ws.Range("L89").Value = 34450
$ws.Range("M89").Value = -54381.855
$ws.Range("N89").Value = -45682
$ws.Range("H94").Value = 712.5
$ws.Range("J94").Value = 300
$ws.Range("L94").Value = 300
$ws.Range("N94").Value = -1202
$ws.Range("H108").Value = 50000
$ws.Range("J108").Value = 50000
$ws.Range("L108").Value = 50000
$ws.Range("N108").Value = -57680
$ws.Range("H132").Value = 2433.2917
$ws.Range("I132").Value = 1847.1765
$ws.Range("K132").Value = 5541.529500000001
$ws.Range("M132").Value = -3011.529500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4999
$ws.Range("I3").Value = 4999
$ws.Range("K3").Value = 14997
$ws.Range("M3").Value = -14885
$ws.Range("H5").Value = 1216.6666
$ws.Range("I5").Value = 825
$ws.Range("K5").Value = 2475
$ws.Range("M5").Value = -2363
$ws.Range("H34").Value = 4113
$ws.Range("I34").Value = 1366
$ws.Range("J34").Value = 6173.25
$ws.Range("K34").Value = 4098
$ws.Range("L34").Value = 18519.75
$ws.Range("M34").Value = -4014
$ws.Range("N34").Value = -18687.75
$ws.Range("H39").Value = 20000
$ws.Range("J39").Value = 20000
$ws.Range("L39").Value = 60000
$ws.Range("N39").Value = -60588
$ws.Range("H55").Value = 1445
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H56").Value = 18478.059
$ws.Range("I56").Value = 18478.059
$ws.Range("K56").Value = 18478.059
$ws.Range("M56").Value = -17948.059
$ws.Range("H107").Value = 196.71428
$ws.Range("J107").Value = 212.83333
$ws.Range("L107").Value = 638.49999
$ws.Range("N107").Value = -4478.49999
$ws.Range("H135").Value = 1216.6666
$ws.Range("I135").Value = 825
$ws.Range("K135").Value = 7425
$ws.Range("M135").Value = -4890

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 23628.666
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 34943
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 34943
$ws.Range("M46").Value = -844
$ws.Range("N46").Value = -35255
$ws.Range("H80").Value = 5358.8184
$ws.Range("I80").Value = 4994.4287
$ws.Range("J80").Value = 5996.5
$ws.Range("K80").Value = 4994.4287
$ws.Range("L80").Value = 5996.5
$ws.Range("M80").Value = -3996.4287
$ws.Range("N80").Value = -7992.5
$ws.Range("H83").Value = 5358.8184
$ws.Range("I83").Value = 4994.4287
$ws.Range("J83").Value = 5996.5
$ws.Range("K83").Value = 24972.1435
$ws.Range("L83").Value = 29982.5
$ws.Range("M83").Value = -19980.1435
$ws.Range("N83").Value = -39966.5
$ws.Range("H97").Value = 550.38464
$ws.Range("I97").Value = 778.25
$ws.Range("J97").Value = 185.8
$ws.Range("K97").Value = 778.25
$ws.Range("L97").Value = 185.8
$ws.Range("M97").Value = -282.25
$ws.Range("N97").Value = -1177.8
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 1943.3334
$ws.Range("I132").Value = 1228.7059
$ws.Range("J132").Value = 3158.2
$ws.Range("K132").Value = 3686.1177
$ws.Range("L132").Value = 9474.599999999999
$ws.Range("M132").Value = -1156.1177
$ws.Range("N132").Value = -14534.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2167.8462
$ws.Range("J46").Value = 2499
$ws.Range("L46").Value = 2499
$ws.Range("N46").Value = -2875
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H93").Value = 4050.2
$ws.Range("I93").Value = 4312.75
$ws.Range("K93").Value = 4312.75
$ws.Range("M93").Value = -3064.75
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 4046.889
$ws.Range("I132").Value = 3996.6667
$ws.Range("K132").Value = 11990.0001
$ws.Range("M132").Value = -9460.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1743.125
$ws.Range("I113").Value = 1708.6666
$ws.Range("K113").Value = 5125.9998
$ws.Range("M113").Value = -2955.9998
$ws.Range("H126").Value = 1842.4286
$ws.Range("I126").Value = 1842.4286
$ws.Range("K126").Value = 5527.2858
$ws.Range("M126").Value = -3057.2858
